# The user cut the bottom block of the "Edit Repayment Schedule" sheet
# (the "approve / disburse" automation steps, rows 8-13) out into a brand
# new "Sheet1" tab appended at the end of the workbook, then left that new
# sheet as the active tab.

$wb = $excel.ActiveWorkbook

$editSheet = $wb.Worksheets.Item("Edit Repayment Schedule")

# Add a new worksheet after the last tab in the workbook -> becomes "Sheet1"
# and is left as the active/selected sheet (matches activeTab + tabSelected
# moving off of NewLoanInput in the diff).
$newSheet = $wb.Worksheets.Add($null, $editSheet)

# Move (cut/paste) A8:B13 from "Edit Repayment Schedule" to A1 on the new
# sheet, then remove the now-empty rows from the source sheet.
$srcRange = $editSheet.Range("A8:B13")
$srcRange.Cut($newSheet.Range("A1")) | Out-Null
$editSheet.Rows("8:13").Delete() | Out-Null

# Restore the recorded selections on both sheets.
$editSheet.Range("A8:XFD18").Select() | Out-Null
$newSheet.Range("A3").Select() | Out-Null
